# [FEATURE] Se agregaron casos de los modulos ASUC 42 - 43 Y CP07
# Adds three new user rows (41-43) to the "Users" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# Row 41: JBORDOY / 173 (numeric password, keep same style as the rows above it)
$ws.Range("C41").Value = 173
$ws.Range("C41").NumberFormat = $ws.Range("C40").NumberFormat()
$ws.Range("C41").HorizontalAlignment = $ws.Range("C40").HorizontalAlignment()
$ws.Range("A41").Value = "JBORDOY"

# Row 42: CRECERAC / Usuario Emergencia (text password)
$ws.Range("C42").NumberFormat = $ws.Range("C40").NumberFormat()
$ws.Range("C42").HorizontalAlignment = $ws.Range("C40").HorizontalAlignment()
$ws.Range("C42").Value = "Usuario Emergencia"
$ws.Range("A42").Value = "CRECERAC"

# Row 43: F00074 / 074 (text password, must keep leading zero)
$ws.Range("C43").NumberFormat = $ws.Range("C40").NumberFormat()
$ws.Range("C43").HorizontalAlignment = $ws.Range("C40").HorizontalAlignment()
$ws.Range("C43").Value = "074"
$ws.Range("A43").Value = "F00074"

# Match the author's final selection/scroll position from the commit
$ws.Range("C42").Select()
